$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new course entry as a new row (row 16)
$ws.Range("A16").Value = "Linear Classifiers in Python"
$ws.Range("B16").Value = 4

# Match the existing "section header" font color used by the rows above
# (rows 13-15, which share the same dark navy blue font color)
$ws.Range("A16").Font.Color = $ws.Range("A15").Font.Color

# The rating cell gets its own (new) font entry, same as black/automatic text
$ws.Range("B16").Font.Color = 0

# Move/update the active selection to the next empty row, as Excel does
# after data entry
$ws.Range("A17").Select() | Out-Null
